# Add 2022-Q1 sheet and update the 总计 (totals) sheet.
#
# Resulting sheet order: 2021-Q4, 2022-Q1, 总计

$wb = $excel.ActiveWorkbook

$q4Sheet    = $wb.Worksheets.Item("2021-Q4")
$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right before the "总计" sheet.
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($totalSheet)
$newSheet.Name = "2022-Q1"

# NOTE: after Worksheets.Add() the old $totalSheet handle now tracks the
# sheet at that same position, which is the newly inserted sheet - re-fetch
# the "总计" sheet by name so later writes land on the right tab.
$totalSheet = $wb.Worksheets.Item("总计")

# Copy header-row (B1:H1) + column-A data style from the "2021-Q4" sheet so the
# new sheet's look matches the existing ones exactly (bold/border/centered style).
$q4Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$q4Sheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Fill in the header row.
# ---------------------------------------------------------------------------
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 3. Fill in the single data row (values are stored as text, matching the
#    source data export format, except the numeric rank in column H).
# ---------------------------------------------------------------------------
$newSheet.Range("A2").Value = 0

$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "165524"

$newSheet.Range("C2").Value = "信诚中证智能家居指数（LOF）"

$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.40"

$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "93.89"

$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "1.32"

$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0053"

$newSheet.Range("H2").Value = 3

# Re-stamp B2:G2 with a pristine, never-touched cell's formatting so that
# forcing text above (via NumberFormat "@") doesn't leave a stray style
# behind. NOTE: only touch formatting here - re-assigning .Value afterwards
# would flip the cell back to "General" and Excel would re-parse these
# numeric-looking strings as numbers again.
$newSheet.Range("Z100").Copy()
$newSheet.Range("B2:G2").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Update the "总计" sheet: push the existing 2021-Q4 total down to row 3
#    and add the new 2022-Q1 total in row 2.
# ---------------------------------------------------------------------------
$totalSheet.Range("B3").Value = "2021-Q4"
$totalSheet.Range("C3").Value = 5
$totalSheet.Range("D3").Value = 0.06
$totalSheet.Range("A3").Value = 1

$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.01
$totalSheet.Range("A2").Value = 0

# Copy the (already correctly styled) A2 cell format down onto the new A3 cell.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$totalSheet.Range("A3").Value = 1
